$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 text changes from "Concentration (aM)" to "Concentration (fM)"
$ws.Range("A1").Value = "Concentration (fM)"
$ws.Range("B1").Value = "Signal"

# Update the selected cell to A2
$ws.Range("A2").Select()
